# Fun Friday game roster update:
#  - akash's meme image was renamed from "akash_deep.png" to "akash.png"
#  - a new employee row for "praveena" was appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update akash's image filename (row 5, column B)
$ws.Range("B5").Value = "meme_images/akash.png"

# Append a new row for praveena
$ws.Range("A25").Value = "praveena"
$ws.Range("B25").Value = "meme_images/praveena.png"

# Match the final active selection left by the author's edit
$ws.Range("H23").Select() | Out-Null
